$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算
# Insert a new row 2 for 2021/12/17 and push the rest of the table down by
# one row. The table also keeps one trailing blank row after the data.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
# B3 already holds the text "202202" (same text value the new row needs);
# copy it into B2 so the new cell keeps the original text formatting
# instead of Excel auto-converting the numeric-looking string to a number.
$ws1.Range("B3").Copy($ws1.Range("B2"))
$ws1.Range("A2").Value = "日期：2021/12/17"
$ws1.Range("C2").Value = 17759
$ws1.Range("D2").Value = 1434
$ws1.Range("E2").Value = 5487531
$ws1.Range("F2").Value = 17645
# Keep the trailing blank row (now row 5) present in the used range.
$ws1.Range("A5").Borders.LineStyle = 0

# ---------------------------------------------------------------------------
# Sheet 2: 散戶多空力道
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").Value = "日期：2021/12/17"
$ws2.Range("B2").Value = 0.17

# ---------------------------------------------------------------------------
# Sheet 3: 三大法人買賣金額
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()
$ws3.Range("A2").Value = "110年12月17日"
$ws3.Range("B2").Value = -25.74
$ws3.Range("C2").Value = 52.59

# ---------------------------------------------------------------------------
# Sheet 4: 大盤多空點位
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
$ws4.Range("A2").Value = "110年12月17日"
$ws4.Range("B2").Value = 17782.21

# ---------------------------------------------------------------------------
# Sheet 5: 期貨大額交易人未沖銷部位
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
# Force text formatting first so the date-like string "2021/12/17" is not
# auto-converted into a date serial number (matches the other rows, which
# all store this column as plain text).
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "2021/12/17"
$ws5.Range("B2").Value = 45843
$ws5.Range("C2").Value = 54309
$ws5.Range("D2").Value = -630
$ws5.Range("E2").Value = -27
$ws5.Range("F2").Value = 22015
$ws5.Range("G2").Value = 47429
$ws5.Range("H2").Value = 812
$ws5.Range("I2").Value = -304
$ws5.Range("J2").Value = -25414
$ws5.Range("K2").Value = 1116
$ws5.Range("L2").Value = -1442
$ws5.Range("M2").Value = 277
$ws5.Range("N2").Value = -1719
